# Apply the commit's changes:
# 1. Update the "Date" metadata value.
# 2. Swap the two "Mapping" columns (AK <-> AL) on the Elements sheet -
#    header text, all data rows, and the column widths - so that
#    "Mapping: Spécification métier..." moves from AL to AK and
#    "Mapping: RIM Mapping" moves from AK to AL.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 : Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value2 = "2024-03-22T16:25:12+00:00"

# --- 2. Elements sheet : swap columns AK (37) and AL (38) ---
$ws = $wb.Worksheets.Item("Elements")

# This sheet has header row 1 plus 20 data rows.
$rowCount = 20

for ($r = 1; $r -le $rowCount; $r++) {
    $akCell = $ws.Cells.Item($r, 37)
    $alCell = $ws.Cells.Item($r, 38)

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    if ($akVal -ne $alVal) {
        $akCell.Value2 = $alVal
        $alCell.Value2 = $akVal
    }
}

# Swap the column widths to match (AK becomes the wide "Spécification" column,
# AL becomes the narrower "RIM Mapping" column).
$ws.Columns.Item(37).ColumnWidth = 85.66666666666667
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668
